$wb = $excel.ActiveWorkbook

# --- TABLE_NAMES sheet: add STUDENT row (creates shared string "STUDENT") ---
$tableNames = $wb.Worksheets.Item("TABLE_NAMES")
$tableNames.Range("A4").Value = "STUDENT"
$tableNames.Range("A4").Select()

# --- EMP sheet: change B2 (empno DATA_TYPE) from varchar to int ---
$emp = $wb.Worksheets.Item("EMP")
$emp.Range("B2").Value = "int"
$emp.Range("C28").Select()

# --- DEPT sheet: update selection ---
$dept = $wb.Worksheets.Item("DEPT")
$dept.Range("B3").Select()

# --- Add STUDENT worksheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$student = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$student.Name = "STUDENT"

# Header row (reuses existing shared strings)
$student.Range("A1").Value = "COULMN_NAME"
$student.Range("B1").Value = "DATA_TYPE"
$student.Range("C1").Value = "DATA_LENGTH"

# Column A values first, in row order, so new shared strings are appended
# in the same order as the target workbook (STUDENT, Student_id, Student_name,
# College, State, Country, then varchar(50)).
$student.Range("A2").Value = "Student_id"
$student.Range("A3").Value = "Student_name"
$student.Range("A4").Value = "College"
$student.Range("A5").Value = "State"
$student.Range("A6").Value = "Country"

$student.Range("B2").Value = "int"
$student.Range("C2").Value = "int"

$student.Range("B3").Value = "varchar"
$student.Range("C3").Value = "varchar(50)"

$student.Range("B4").Value = "varchar"
$student.Range("C4").Value = "varchar(50)"

$student.Range("B5").Value = "varchar"
$student.Range("C5").Value = "varchar(50)"

$student.Range("B6").Value = "varchar"
$student.Range("C6").Value = "varchar(50)"

# Match the wrap-text / vertically-centered formatting used by the other sheets
$body = $student.Range("A2:C6")
$body.VerticalAlignment = -4108
$body.WrapText = $true

# Column widths (approximating the bestFit widths Excel computed for this data)
$student.Columns.Item(1).ColumnWidth = 14.053385416666666
$student.Columns.Item(2).ColumnWidth = 9.608072916666666
$student.Columns.Item(3).ColumnWidth = 12.385416666666666

$student.Range("D12").Select()
